$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reset number format for cells that were scientific-notation (s=4) styled but need General (s=0) ---
$ws.Range("E108").NumberFormat = "GENERAL"
$ws.Range("E129").NumberFormat = "GENERAL"
$ws.Range("E131").NumberFormat = "GENERAL"
$ws.Range("E135").NumberFormat = "GENERAL"
$ws.Range("E136").NumberFormat = "GENERAL"
$ws.Range("E145").NumberFormat = "GENERAL"
$ws.Range("E147").NumberFormat = "GENERAL"
$ws.Range("E148").NumberFormat = "GENERAL"
$ws.Range("E151").NumberFormat = "GENERAL"
$ws.Range("E152").NumberFormat = "GENERAL"
$ws.Range("E156").NumberFormat = "GENERAL"
$ws.Range("E157").NumberFormat = "GENERAL"

# Row 108
$ws.Range("E108").Value = 7.74

# Row 129
$ws.Range("E129").Value = 0.059
$ws.Range("F129").Value = 0.055
$ws.Range("G129").Value = 0.52

# Row 130
$ws.Range("E130").Value = 0.036
$ws.Range("F130").Value = 0.122
$ws.Range("G130").Value = 0.89

# Row 131
$ws.Range("E131").Value = 0.04
$ws.Range("F131").Value = 0.107
$ws.Range("G131").Value = 0.45

# Row 132
$ws.Range("E132").Value = 0.031
$ws.Range("F132").Value = 0.128
$ws.Range("G132").Value = 0.85

# Row 133
$ws.Range("E133").Value = 0.044
$ws.Range("F133").Value = 0.129
$ws.Range("G133").Value = 0.9

# Row 134
$ws.Range("E134").Value = 0.032
$ws.Range("F134").Value = 0.129
$ws.Range("G134").Value = 0.88

# Row 135
$ws.Range("E135").Value = 0.028
$ws.Range("F135").Value = 0.127
$ws.Range("G135").Value = 0.75

# Row 136
$ws.Range("E136").Value = 0.026
$ws.Range("F136").Value = 0.126
$ws.Range("G136").Value = 0.75

# Row 137
$ws.Range("E137").Value = 0.026
$ws.Range("F137").Value = 0.124
$ws.Range("G137").Value = 0.77

# Row 138
$ws.Range("E138").Value = 0.026
$ws.Range("F138").Value = 0.125
$ws.Range("G138").Value = 0.71

# Row 139
$ws.Range("E139").Value = 0.027
$ws.Range("F139").Value = 0.127
$ws.Range("G139").Value = 0.78

# Row 140
$ws.Range("E140").Value = 0.025
$ws.Range("F140").Value = 0.128
$ws.Range("G140").Value = 0.74

# Row 141
$ws.Range("E141").Value = 0.024
$ws.Range("F141").Value = 0.123
$ws.Range("G141").Value = 0.63

# Row 142
$ws.Range("E142").Value = 0.023
$ws.Range("F142").Value = 0.122
$ws.Range("G142").Value = 0.71

# Row 143
$ws.Range("E143").Value = 0.043
$ws.Range("F143").Value = 0.015
$ws.Range("G143").Value = 0.16
$ws.Range("I143").Value = "converges to same point"

# Row 144
$ws.Range("E144").Value = 0.043
$ws.Range("F144").Value = 0.015
$ws.Range("G144").Value = 0.13
$ws.Range("I144").Value = "converges to same point"

# Row 145
$ws.Range("E145").Value = 0.022
$ws.Range("F145").Value = 0.136
$ws.Range("G145").Value = 0.66

# Row 146
$ws.Range("E146").Value = 0.036
$ws.Range("F146").Value = 0.061
$ws.Range("G146").Value = 0.2

# Row 147
$ws.Range("E147").Value = 0.028
$ws.Range("F147").Value = 0.12
$ws.Range("G147").Value = 0.52

# Row 148
$ws.Range("E148").Value = 0.043
$ws.Range("F148").Value = 0.015
$ws.Range("G148").Value = 0.09
$ws.Range("I148").Value = "converges to same point"

# Row 149
$ws.Range("E149").Value = 0.043
$ws.Range("F149").Value = 0.015
$ws.Range("G149").Value = 0.08
$ws.Range("I149").Value = "converges to same point"

# Row 150
$ws.Range("E150").Value = 0.035
$ws.Range("F150").Value = 0.068
$ws.Range("G150").Value = 0.11

# Row 151
$ws.Range("E151").Value = 0.042
$ws.Range("F151").Value = 0.015
$ws.Range("G151").Value = 0.12
$ws.Range("I151").Value = "converges to same point"

# Row 152
$ws.Range("E152").Value = 0.024
$ws.Range("F152").Value = 0.127
$ws.Range("G152").Value = 0.44

# Row 153
$ws.Range("E153").Value = 0.129
$ws.Range("F153").Value = 0.04
$ws.Range("G153").Value = 0.05
$ws.Range("I153").Value = "converges to same point"

# Row 154
$ws.Range("E154").Value = 0.042
$ws.Range("F154").Value = 0.015
$ws.Range("G154").Value = 0.11
$ws.Range("I154").Value = "converges to same point"

# Row 155
$ws.Range("E155").Value = 162.36
$ws.Range("F155").Value = 0.062
$ws.Range("G155").Value = 0.95
$ws.Range("H155").Value = "higher"
$ws.Range("I155").Value = "DIVERGED"

# Row 156
$ws.Range("C156").Value = "run112"
$ws.Range("E156").Value = 0.042
$ws.Range("F156").Value = 0.015
$ws.Range("G156").Value = 0.14
$ws.Range("I156").Value = "converges to same point"

# Row 157
$ws.Range("C157").Value = "run113"
$ws.Range("E157").Value = 2443.44
$ws.Range("F157").Value = 0.085
$ws.Range("G157").Value = 0.95
$ws.Range("H157").Value = "higher"
$ws.Range("I157").Value = "DIVERGED"

# Row 158
$ws.Range("C158").Value = "run114"
$ws.Range("E158").Value = 174586
$ws.Range("F158").Value = 0.058
$ws.Range("G158").Value = 0.95
$ws.Range("H158").Value = "higher"
$ws.Range("I158").Value = "DIVERGED"

# Row 160
$ws.Range("G160").Value = "Notes:"

# Row 161
$ws.Range("B161").Value = "10 .^ unifrnd(-6,-2,30,1)"
$ws.Range("G161").Value = "For plots deleted points with 6 highest alphas or 6 highest lambdas (if repreated)"

# Row 162
$ws.Range("B162").Value = "10 .^ unifrnd(-5, 0, 30, 1)"

# Row 166
$ws.Range("A166").Value = 0.00004
$ws.Range("B166").Value = 0.00007
$ws.Range("C166").Value = "run115"
$ws.Range("D166").Value = 1

# --- Update selection / scroll position (cosmetic) ---
$ws.Range("D166").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 131
$win.ScrollColumn = 1
